$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Row 2
$ws.Range("B2").Value = "2025-11-27T17:00:00"
$ws.Range("C2").Value = "Салават Юлаев"
$ws.Range("D2").Value = "Барыс"
$ws.Range("E2").Value = 897819
$ws.Range("F2").Value = "https://text.khl.ru/text/897819.html"
$ws.Range("G2").Value = 1.454545
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 2.59603
$ws.Range("J2").Value = 4.227273
$ws.Range("K2").Value = 2.840909
$ws.Range("L2").Value = 1.798015
$ws.Range("M2").Value = 2.454545
$ws.Range("N2").Value = 23.48092
$ws.Range("O2").Value = 23.645876
$ws.Range("P2").Value = 47.126796
$ws.Range("Q2").Value = -0.2
$ws.Range("R2").Value = -0.2
$ws.Range("S2").Value = 0.59579
$ws.Range("T2").Value = 0.172039
$ws.Range("U2").Value = 0.231982
$ws.Range("V2").Value = 0.319402
$ws.Range("W2").Value = 0.68041
$ws.Range("X2").Value = 0.505953
$ws.Range("Y2").Value = 0.493858
$ws.Range("Z2").Value = 0.679033
$ws.Range("AA2").Value = 0.320779
$ws.Range("AB2").Value = 0.81285
$ws.Range("AC2").Value = 0.186961
$ws.Range("AD2").Value = 0.901531
$ws.Range("AE2").Value = 0.09828
$ws.Range("AF2").Value = 0.775796
$ws.Range("AG2").Value = 0.224204
$ws.Range("AH2").Value = 0.540241
$ws.Range("AI2").Value = 0.459759
$ws.Range("AJ2").Value = 0.536572
$ws.Range("AK2").Value = 0.463428
$ws.Range("AL2").Value = 0.268847
$ws.Range("AM2").Value = 0.731153
$ws.Range("AN2").Value = 0.888398
$ws.Range("AO2").Value = 0.594522

# Row 3
$ws.Range("B3").Value = "2025-11-27T17:00:00"
$ws.Range("C3").Value = "Металлург Мг"
$ws.Range("D3").Value = "Авангард"
$ws.Range("E3").Value = 897821
$ws.Range("F3").Value = "https://text.khl.ru/text/897821.html"
$ws.Range("G3").Value = 6.3
$ws.Range("H3").Value = 5.678571
$ws.Range("I3").Value = 4.6
$ws.Range("J3").Value = 1.957785
$ws.Range("K3").Value = 4.128893
$ws.Range("L3").Value = 5.139286
$ws.Range("M3").Value = 11.978571
$ws.Range("N3").Value = 37.406257
$ws.Range("O3").Value = 35.92968
$ws.Range("P3").Value = 73.335938
$ws.Range("Q3").Value = 0.2
$ws.Range("R3").Value = 0.2
$ws.Range("S3").Value = 0.304904
$ws.Range("T3").Value = 0.126149
$ws.Range("U3").Value = 0.549013
$ws.Range("V3").Value = 0.017546
$ws.Range("W3").Value = 0.96252
$ws.Range("X3").Value = 0.046563
$ws.Range("Y3").Value = 0.933504
$ws.Range("Z3").Value = 0.100349
$ws.Range("AA3").Value = 0.879717
$ws.Range("AB3").Value = 0.183432
$ws.Range("AC3").Value = 0.796634
$ws.Range("AD3").Value = 0.293437
$ws.Range("AE3").Value = 0.686629
$ws.Range("AF3").Value = 0.917421
$ws.Range("AG3").Value = 0.082579
$ws.Range("AH3").Value = 0.780181
$ws.Range("AI3").Value = 0.219819
$ws.Range("AJ3").Value = 0.964012
$ws.Range("AK3").Value = 0.035988
$ws.Range("AL3").Value = 0.8866
$ws.Range("AM3").Value = 0.1134
$ws.Range("AN3").Value = 0.563876
$ws.Range("AO3").Value = 0.781871

# Row 4
$ws.Range("B4").Value = "2025-11-27T19:30:00"
$ws.Range("C4").Value = "Спартак"
$ws.Range("D4").Value = "Автомобилист"
$ws.Range("E4").Value = 897820
$ws.Range("F4").Value = "https://text.khl.ru/text/897820.html"
$ws.Range("G4").Value = 4.07223
$ws.Range("H4").Value = 2.501555
$ws.Range("I4").Value = 4.583333
$ws.Range("J4").Value = 2.83482
$ws.Range("K4").Value = 3.453525
$ws.Range("L4").Value = 3.542444
$ws.Range("M4").Value = 6.573785
$ws.Range("N4").Value = 36.594732
$ws.Range("O4").Value = 28.024318
$ws.Range("P4").Value = 64.61905
$ws.Range("Q4").Value = 0.153209
$ws.Range("R4").Value = -0.038073
$ws.Range("S4").Value = 0.409069
$ws.Range("T4").Value = 0.153703
$ws.Range("U4").Value = 0.435191
$ws.Range("V4").Value = 0.081976
$ws.Range("W4").Value = 0.915987
$ws.Range("X4").Value = 0.17336
$ws.Range("Y4").Value = 0.824603
$ws.Range("Z4").Value = 0.301223
$ws.Range("AA4").Value = 0.69674
$ws.Range("AB4").Value = 0.450312
$ws.Range("AC4").Value = 0.547651
$ws.Range("AD4").Value = 0.599314
$ws.Range("AE4").Value = 0.398649
$ws.Range("AF4").Value = 0.859117
$ws.Range("AG4").Value = 0.140883
$ws.Range("AH4").Value = 0.670471
$ws.Range("AI4").Value = 0.329529
$ws.Range("AJ4").Value = 0.86853
$ws.Range("AK4").Value = 0.13147
$ws.Range("AL4").Value = 0.686932
$ws.Range("AM4").Value = 0.313068
$ws.Range("AN4").Value = 0.706839
$ws.Range("AO4").Value = 0.729345

# Row 5
$ws.Range("B5").Value = "2025-11-27T19:30:00"
$ws.Range("C5").Value = "ЦСКА"
$ws.Range("D5").Value = "Лада"
$ws.Range("E5").Value = 897822
$ws.Range("F5").Value = "https://text.khl.ru/text/897822.html"
$ws.Range("G5").Value = 2.188722
$ws.Range("H5").Value = 1.117647
$ws.Range("I5").Value = 1.178571
$ws.Range("J5").Value = 3.442604
$ws.Range("K5").Value = 2.815663
$ws.Range("L5").Value = 1.148109
$ws.Range("M5").Value = 3.306369
$ws.Range("N5").Value = 27.591019
$ws.Range("O5").Value = 20.850574
$ws.Range("P5").Value = 48.441593
$ws.Range("Q5").Value = -0.101053
$ws.Range("R5").Value = -0.2
$ws.Range("S5").Value = 0.719666
$ws.Range("T5").Value = 0.151933
$ws.Range("U5").Value = 0.128229
$ws.Range("V5").Value = 0.44058
$ws.Range("W5").Value = 0.559249
$ws.Range("X5").Value = 0.635914
$ws.Range("Y5").Value = 0.363914
$ws.Range("Z5").Value = 0.790767
$ws.Range("AA5").Value = 0.209062
$ws.Range("AB5").Value = 0.893066
$ws.Range("AC5").Value = 0.106762
$ws.Range("AD5").Value = 0.950994
$ws.Range("AE5").Value = 0.048834
$ws.Range("AF5").Value = 0.771575
$ws.Range("AG5").Value = 0.228425
$ws.Range("AH5").Value = 0.534272
$ws.Range("AI5").Value = 0.465728
$ws.Range("AJ5").Value = 0.318542
$ws.Range("AK5").Value = 0.681458
$ws.Range("AL5").Value = 0.109459
$ws.Range("AM5").Value = 0.890541
$ws.Range("AN5").Value = 0.953673
$ws.Range("AO5").Value = 0.481443
